$wb = $excel.ActiveWorkbook
$data = $wb.Worksheets.Item("data")

# --- Update time_taken values for the "data" sheet (column F), rows 2-19 ---
$newTimes = @(
    "2021-10-05 14:20:56.067227",
    "2021-10-05 14:20:56.067234",
    "2021-10-05 14:20:56.067238",
    "2021-10-05 14:20:56.067240",
    "2021-10-05 14:20:56.067243",
    "2021-10-05 14:20:56.067246",
    "2021-10-05 14:20:56.067248",
    "2021-10-05 14:20:56.067251",
    "2021-10-05 14:20:56.067253",
    "2021-10-05 14:20:56.067256",
    "2021-10-05 14:20:56.067259",
    "2021-10-05 14:20:56.067261",
    "2021-10-05 14:20:56.067264",
    "2021-10-05 14:20:56.067266",
    "2021-10-05 14:20:56.067269",
    "2021-10-05 14:20:56.067271",
    "2021-10-05 14:20:56.067274",
    "2021-10-05 14:20:56.067276"
)

for ($i = 0; $i -lt $newTimes.Length; $i++) {
    $row = $i + 2
    $data.Cells.Item($row, 6).Value = $newTimes[$i]
}

# --- Add the new "metadata" worksheet, placed right after "data" ---
$meta = $wb.Worksheets.Add($null, $data)
$meta.Name = "metadata"

# Header row (B1:G1)
$meta.Cells.Item(1, 2).Value = "data_name"
$meta.Cells.Item(1, 3).Value = "data_id"
$meta.Cells.Item(1, 4).Value = "data_version"
$meta.Cells.Item(1, 5).Value = "data_version_created"
$meta.Cells.Item(1, 6).Value = "panel_query_time"
$meta.Cells.Item(1, 7).Value = "panel_get_request"

# Data row (row 2)
$meta.Cells.Item(2, 1).Value = 0
$meta.Cells.Item(2, 2).Value = "Hypophosphataemia or rickets"
$meta.Cells.Item(2, 3).Value = 482

# data_version must stay textual ("2.14"), not be coerced to a number
$meta.Cells.Item(2, 4).NumberFormat = "@"
$meta.Cells.Item(2, 4).Value = "2.14"

$meta.Cells.Item(2, 5).Value = "2020-11-30T10:33:04.911127Z"
$meta.Cells.Item(2, 6).Value = "2021-10-05 14:20:56.063616"
$meta.Cells.Item(2, 7).Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/482/?format=json"

# Copy the header/index styling from the "data" sheet (bold, centered, bordered)
$data.Range("B1").Copy()
$meta.Range("B1:G1").PasteSpecial(-4122)  # xlPasteFormats

$data.Range("A2").Copy()
$meta.Range("A2").PasteSpecial(-4122)     # xlPasteFormats

# Keep "data" as the active/selected sheet (unchanged from the original workbook)
$data.Activate()
